$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing A21 timestamp value
$ws.Range("A21").Value = 45874.83353430556

# Add the new row 22 data
$ws.Range("A22").Value = 45874.87520286583
$ws.Range("A22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B22").Value = 2025
$ws.Range("C22").Value = 19
$ws.Range("D22").Value = 15
$ws.Range("E22").Value = 86.68000000000001
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 1.42
$ws.Range("H22").Value = "E"
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = "21:00:17"
